{"js": "// Update questions document with performance metrics for different tested algorithms.\n\n// ---------------------------------------------------------------------\n// Change 1: paragraph about `SelectKBest` / normalizing features.\n// Add the Gaussian NB precision/recall callouts around the k=10 and k=5\n// sentences.\n// ---------------------------------------------------------------------\n{\n  const body = context.document.body;\n\n  const r1 = body.search(\n    \"Initially I had it select 10 features, but the performance of the algorithm was not sufficient. When I had\",\n    { matchCase: true }\n  );\n  r1.load(\"items\");\n  await context.sync();\n  if (r1.items.length > 0) {\n    r1.items[0].insertText(\n      \"Initially I had it select 10 features while using a Gaussian Naive Bayes classifier, but the performance of the algorithm was not sufficient (precision: 0.33357, recall: 0.23100). When I had\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n\n  const r2 = body.search(\n    \"When I reduced k to 5, I saw improved recall and precision. The top five features\",\n    { matchCase: true }\n  );\n  r2.load(\"items\");\n  await context.sync();\n  if (r2.items.length > 0) {\n    r2.items[0].insertText(\n      \"When I reduced K to 5, I saw improved recall and precision for the Gaussian Naive Bayes classifier (precision: 0.47400, recall: 35100). The top five features\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// Change 2: paragraph \"I ended up using a Decision Tree...\" is rewritten\n// and split into three paragraphs: the classifier-comparison paragraph\n// (Gaussian NB / Decision Tree numbers), a new SVM paragraph, and a new\n// paragraph about the `kfold_eval` evaluation.\n// ---------------------------------------------------------------------\n{\n  const body = context.document.body;\n\n  const r3 = body.search(\n    \"I ended up using a Decision Tree with the `min_samples_split` parameter set to 5. I tried a Decision Tree with different values for that parameter as recommended by my parameter tuning with `GridSearchCV`. I re-tuned the algorithm a few times during the project because my feature selection (as described in the previous answer) was incremental. I also tried a Gaussian Naive Bayes classifier and a Support Vector Machine classifier. They tended to have better accuracy scores (particularly the SVM), but their precision and recall were not better than the Decision Tree classifier.\",\n    { matchCase: true }\n  );\n  r3.load(\"items\");\n  await context.sync();\n  if (r3.items.length > 0) {\n    r3.items[0].insertText(\n      \"I ended up using a Gaussian Naive Bayes classifier. I also tried Decision Tree classifier with different values for `min_samples_split` as recommended by my parameter tuning with `GridSearchCV`. With the 5 best features (as described in the previous answer) and `min_samples_split` equal to 2, precision was 0.27428 and recall was 0.27250. Setting `min_samples_split` to 4 yield a precision of 0.266601 and recall of 0.25750; `min_samples_split` set to 5 gave a precision of 0.27428 and recall of 27250. These were worse than my Gaussian Naive Bayes algorithm.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n\n  // Remove the trailing \"I evaluated this with a custom function...\" run\n  // (its content moves into the new third paragraph below).\n  const r4 = body.search(\n    \"I evaluated this with a custom function `kfold_eval`, which deployed K-fold cross-validation and reported average scores, precision, and recall.\",\n    { matchCase: true }\n  );\n  r4.load(\"items\");\n  await context.sync();\n  if (r4.items.length > 0) {\n    r4.items[0].delete();\n    await context.sync();\n  }\n\n  // Find the (now-edited) paragraph again and insert the two new\n  // paragraphs after it, matching the inherited TextBody/ind formatting.\n  const paras = body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n\n  let targetIndex = -1;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.indexOf(\"I ended up using a Gaussian Naive Bayes classifier.\") === 0) {\n      targetIndex = i;\n      break;\n    }\n  }\n\n  if (targetIndex !== -1) {\n    const anchorPara = paras.items[targetIndex];\n    const svmPara = anchorPara.insertParagraph(\n      \"I also tried a a Support Vector Machine classifier. Out of the box it tended to have better accuracy scores, but its precision and recall were not better. Trying to tune it with `GridSearchCV` was painfully slow so I opted not to continue with an SVM. \",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n\n    svmPara.insertParagraph(\n      \"Before passing my algorithm and features to `tester.py`, I evaluated my algorithm and features with a custom function `kfold_eval`, which deployed K-fold cross-validation and reported average scores, precision, and recall.\",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# Update questions document with performance metrics for different tested algorithms.\n\n$d = $word.ActiveDocument\n\nfunction Find-And-Replace($doc, $searchText, $replaceText) {\n    $rng = $doc.Content\n    $find = $rng.Find\n    $find.Text = $searchText\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Text = $replaceText\n    }\n    return $found\n}\n\n# ---------------------------------------------------------------------\n# Change 1: paragraph about `SelectKBest` / normalizing features.\n# Add the Gaussian NB precision/recall callouts around the k=10 and k=5\n# sentences.\n# ---------------------------------------------------------------------\n\nFind-And-Replace $d 'Initially I had it select 10 features, but the performance of the algorithm was not sufficient. When I had' 'Initially I had it select 10 features while using a Gaussian Naive Bayes classifier, but the performance of the algorithm was not sufficient (precision: 0.33357, recall: 0.23100). When I had' | Out-Null\n\nFind-And-Replace $d 'When I reduced k to 5, I saw improved recall and precision. The top five features' 'When I reduced K to 5, I saw improved recall and precision for the Gaussian Naive Bayes classifier (precision: 0.47400, recall: 35100). The top five features' | Out-Null\n\n# ---------------------------------------------------------------------\n# Change 2: paragraph \"I ended up using a Decision Tree...\" is rewritten\n# and split into three paragraphs: the classifier-comparison paragraph\n# (Gaussian NB / Decision Tree numbers), a new SVM paragraph, and a new\n# paragraph about the `kfold_eval` evaluation.\n# ---------------------------------------------------------------------\n\n$old2 = 'I ended up using a Decision Tree with the `min_samples_split` parameter set to 5. I tried a Decision Tree with different values for that parameter as recommended by my parameter tuning with `GridSearchCV`. I re-tuned the algorithm a few times during the project because my feature selection (as described in the previous answer) was incremental. I also tried a Gaussian Naive Bayes classifier and a Support Vector Machine classifier. They tended to have better accuracy scores (particularly the SVM), but their precision and recall were not better than the Decision Tree classifier. '\n$new2 = 'I ended up using a Gaussian Naive Bayes classifier. I also tried Decision Tree classifier with different values for `min_samples_split` as recommended by my parameter tuning with `GridSearchCV`. With the 5 best features (as described in the previous answer) and `min_samples_split` equal to 2, precision was 0.27428 and recall was 0.27250. Setting `min_samples_split` to 4 yield a precision of 0.266601 and recall of 0.25750; `min_samples_split` set to 5 gave a precision of 0.27428 and recall of 27250. These were worse than my Gaussian Naive Bayes algorithm. '\nFind-And-Replace $d $old2 $new2 | Out-Null\n\n# Remove the trailing \"I evaluated this with a custom function...\" sentence\n# (its content moves into the new third paragraph below).\n$oldTail = 'I evaluated this with a custom function `kfold_eval`, which deployed K-fold cross-validation and reported average scores, precision, and recall.'\nFind-And-Replace $d $oldTail '' | Out-Null\n\n# Locate the (now-edited) paragraph again and insert the two new\n# paragraphs after it, inheriting the TextBody/ind formatting.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.IndexOf('I ended up using a Gaussian Naive Bayes classifier.') -eq 0) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ne -1) {\n    $anchorPara = $d.Paragraphs.Item($targetIndex)\n    $anchorPara.Range.InsertParagraphAfter()\n    $svmPara = $d.Paragraphs.Item($targetIndex + 1)\n    $svmPara.Range.Text = 'I also tried a a Support Vector Machine classifier. Out of the box it tended to have better accuracy scores, but its precision and recall were not better. Trying to tune it with `GridSearchCV` was painfully slow so I opted not to continue with an SVM. '\n\n    $svmPara.Range.InsertParagraphAfter()\n    $kfoldPara = $d.Paragraphs.Item($targetIndex + 2)\n    $kfoldPara.Range.Text = 'Before passing my algorithm and features to `tester.py`, I evaluated my algorithm and features with a custom function `kfold_eval`, which deployed K-fold cross-validation and reported average scores, precision, and recall.'\n}\n"}
